{"js": "// Delete the paragraph containing \"hola\" (a leftover placeholder paragraph).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"hola\") {\n    paragraphs.items[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the leftover \"hola\" paragraph (placeholder text left in the template).\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$searchRange.Find.ClearFormatting()\n$found = $searchRange.Find.Execute(\"hola\")\n\nif ($searchRange.Find.Found) {\n    # Paragraphs(1) is the paragraph containing the found text; deleting its\n    # Range removes the run and the paragraph mark, merging it away entirely.\n    $para = $searchRange.Paragraphs(1)\n    $para.Range.Delete()\n}\n"}
